$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.372.48"
$ws.Range("E2").Value = "  +3.43%  "
Set-TextValue $ws.Range("D3") "3.619.10"
$ws.Range("E3").Value = "  +2.58%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "627.59"
$ws.Range("E5").Value = "  +3.22%  "
Set-TextValue $ws.Range("D6") "159.39"
$ws.Range("E6").Value = "  +4.55%  "
Set-TextValue $ws.Range("D7") "3.619.04"
$ws.Range("E7").Value = "  +2.68%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  +2.08%  "
Set-TextValue $ws.Range("D10") "0.146"
$ws.Range("E10").Value = "  +4.53%  "
Set-TextValue $ws.Range("D11") "7.24"
$ws.Range("E11").Value = "  +6.57%  "
$ws.Range("E12").Value = "  +3.28%  "
Set-TextValue $ws.Range("D13") "0.0000224"
$ws.Range("E13").Value = "  +1.59%  "
Set-TextValue $ws.Range("D14") "33.19"
$ws.Range("E14").Value = "  +5.02%  "
Set-TextValue $ws.Range("D15") "4.231.60"
$ws.Range("E15").Value = "  +2.52%  "
Set-TextValue $ws.Range("D16") "3.605.70"
$ws.Range("E16").Value = "  +1.84%  "
Set-TextValue $ws.Range("D17") "69.069.00"
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("E18").Value = "  -0.31%  "
Set-TextValue $ws.Range("D19") "6.62"
$ws.Range("E19").Value = "  +5.27%  "
Set-TextValue $ws.Range("D20") "15.94"
$ws.Range("E20").Value = "  +3.82%  "
Set-TextValue $ws.Range("D21") "10.10"
$ws.Range("E21").Value = "  +9.60%  "
Set-TextValue $ws.Range("D22") "459.90"
$ws.Range("E22").Value = "  +4.04%  "
Set-TextValue $ws.Range("D23") "0.639"
$ws.Range("E23").Value = "  +2.13%  "
Set-TextValue $ws.Range("D24") "78.50"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("E25").Value = "  +12.84%  "
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue $ws.Range("D26") "3.763.50"
$ws.Range("E26").Value = "  +2.46%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D27") "10.62"
$ws.Range("E27").Value = "  +4.65%  "
$ws.Range("E28").Value = "  +0.03%  "
Set-TextValue $ws.Range("D29") "9.24"
$ws.Range("E29").Value = "  +12.80%  "
$ws.Range("E30").Value = "  +3.26%  "
Set-TextValue $ws.Range("D31") "1.71"
$ws.Range("E31").Value = "  +3.04%  "
$ws.Range("E32").Value = "  +11.83%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D33") "1.01"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D34") "6.57"
$ws.Range("E34").Value = "  +7.50%  "
$ws.Range("E35").Value = "  +5.77%  "
Set-TextValue $ws.Range("D36") "26.46"
$ws.Range("E36").Value = "  +3.27%  "
Set-TextValue $ws.Range("D37") "3.605.01"
$ws.Range("E37").Value = "  +2.27%  "
Set-TextValue $ws.Range("D38") "8.37"
$ws.Range("E38").Value = "  +5.30%  "
Set-TextValue $ws.Range("D39") "2.38"
$ws.Range("E39").Value = "  +11.16%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  +7.85%  "
Set-TextValue $ws.Range("D42") "1.00"
$ws.Range("E42").Value = "  +0.05%  "
Set-TextValue $ws.Range("D43") "176.77"
$ws.Range("E43").Value = "  +0.75%  "
Set-TextValue $ws.Range("D44") "5.61"
$ws.Range("E44").Value = "  +1.47%  "
Set-TextValue $ws.Range("D45") "31.86"
$ws.Range("E45").Value = "  +16.21%  "
Set-TextValue $ws.Range("D46") "0.913"
$ws.Range("E46").Value = "  +2.72%  "
Set-TextValue $ws.Range("D47") "1.37"
$ws.Range("E47").Value = "  +12.81%  "
Set-TextValue $ws.Range("D48") "2.80"
$ws.Range("E48").Value = "  +8.58%  "
Set-TextValue $ws.Range("D49") "46.31"
$ws.Range("E49").Value = "  +1.75%  "
Set-TextValue $ws.Range("D50") "7.79"
$ws.Range("E50").Value = "  +3.38%  "
Set-TextValue $ws.Range("D51") "0.266"
$ws.Range("E51").Value = "  +7.61%  "
